$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.032832748375253
$ws.Range("D2").Value = 1.046344668908361
$ws.Range("E2").Value = 1.03195827190779
$ws.Range("F2").Value = 1.053196334805676
$ws.Range("I2").Value = 1.03766955874605
$ws.Range("J2").Value = 1.0379607991498
$ws.Range("K2").Value = 1.049110366666996
$ws.Range("L2").Value = 1.034764887323942
$ws.Range("M2").Value = 1.055942971156751
$ws.Range("N2").Value = 1.03943482254036
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.035573872792656
$ws.Range("D3").Value = 1.047226457215766
$ws.Range("E3").Value = 1.034356095755814
$ws.Range("F3").Value = 1.054590940309874
$ws.Range("I3").Value = 1.037907034404308
$ws.Range("J3").Value = 1.040335306844945
$ws.Range("K3").Value = 1.049803596280695
$ws.Range("L3").Value = 1.03696713436302
$ws.Range("M3").Value = 1.05714908345357
$ws.Range("N3").Value = 1.041812702308792
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.037337268428877
$ws.Range("D4").Value = 1.047793676142792
$ws.Range("E4").Value = 1.035898217808485
$ws.Range("F4").Value = 1.055487825837564
$ws.Range("I4").Value = 1.03805792104658
$ws.Range("J4").Value = 1.041861729071215
$ws.Range("K4").Value = 1.050248177205821
$ws.Range("L4").Value = 1.038382439930384
$ws.Range("M4").Value = 1.057923486629055
$ws.Range("N4").Value = 1.043341292229706
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.038076203230033
$ws.Range("D5").Value = 1.048031343421471
$ws.Range("E5").Value = 1.036544326586038
$ws.Range("F5").Value = 1.055863580229932
$ws.Range("I5").Value = 1.038120694454123
$ws.Range("J5").Value = 1.042501091912006
$ws.Range("K5").Value = 1.050434136969881
$ws.Range("L5").Value = 1.038975168557306
$ws.Range("M5").Value = 1.058247625241122
$ws.Range("N5").Value = 1.043981563039067
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.038200135044488
$ws.Range("D6").Value = 1.04807120273447
$ws.Range("E6").Value = 1.036652683800114
$ws.Range("F6").Value = 1.055926595738892
$ws.Range("I6").Value = 1.038131195872133
$ws.Range("J6").Value = 1.042608307935308
$ws.Range("K6").Value = 1.050465305533412
$ws.Range("L6").Value = 1.039074559022948
$ws.Range("M6").Value = 1.05830196692168
$ws.Range("N6").Value = 1.04408893132142
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.037347151439189
$ws.Range("D7").Value = 1.047796854958576
$ws.Range("E7").Value = 1.035906859710073
$ws.Range("F7").Value = 1.05549285174676
$ws.Range("I7").Value = 1.03805876241162
$ws.Range("J7").Value = 1.041870281404624
$ws.Range("K7").Value = 1.050250665694893
$ws.Range("L7").Value = 1.038390368829078
$ws.Range("M7").Value = 1.057927823335039
$ws.Range("N7").Value = 1.04334985670841
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.033761297189624
$ws.Range("D8").Value = 1.046643375884356
$ws.Range("E8").Value = 1.032770615748698
$ws.Range("F8").Value = 1.053668805666905
$ws.Range("I8").Value = 1.037750392733979
$ws.Range("J8").Value = 1.038765390637913
$ws.Range("K8").Value = 1.049345479256977
$ws.Range("L8").Value = 1.035511188154419
$ws.Range("M8").Value = 1.056351845813849
$ws.Range("N8").Value = 1.040240556640674
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.027360613827603
$ws.Range("D9").Value = 1.044584531170682
$ws.Range("E9").Value = 1.027169267400815
$ws.Range("F9").Value = 1.050411220840563
$ws.Range("I9").Value = 1.037185504894119
$ws.Range("J9").Value = 1.033214553905636
$ws.Range("K9").Value = 1.047719357076175
$ws.Range("L9").Value = 1.030360923846091
$ws.Range("M9").Value = 1.053527493077288
$ws.Range("N9").Value = 1.034681837083548
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.023033987049522
$ws.Range("D10").Value = 1.043193512659956
$ws.Range("E10").Value = 1.023380906019408
$ws.Range("F10").Value = 1.048208727077987
$ws.Range("I10").Value = 1.03679411318609
$ws.Range("J10").Value = 1.029456619023988
$ws.Range("K10").Value = 1.046613605020734
$ws.Range("L10").Value = 1.026872205417472
$ws.Range("M10").Value = 1.051611265008568
$ws.Range("N10").Value = 1.030918565503327
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.021145394948956
$ws.Range("D11").Value = 1.04258662687459
$ws.Range("E11").Value = 1.021726805939236
$ws.Range("F11").Value = 1.047247359237911
$ws.Range("I11").Value = 1.036621047130666
$ws.Range("J11").Value = 1.027814909035591
$ws.Range("K11").Value = 1.046129483233999
$ws.Range("L11").Value = 1.025347641835661
$ws.Range("M11").Value = 1.050773261709678
$ws.Range("N11").Value = 1.029274524098439
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.020441526476554
$ws.Range("D12").Value = 1.042360500489943
$ws.Range("E12").Value = 1.021110263936215
$ws.Range("F12").Value = 1.046889079491171
$ws.Range("I12").Value = 1.036556216234645
$ws.Range("J12").Value = 1.027202850068592
$ws.Range("K12").Value = 1.045948843227012
$ws.Range("L12").Value = 1.024779187591816
$ws.Range("M12").Value = 1.050460717608458
$ws.Range("N12").Value = 1.028661595937503
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.020592616937382
$ws.Range("D13").Value = 1.042409037396158
$ws.Range("E13").Value = 1.021242612156682
$ws.Range("F13").Value = 1.046965985858051
$ws.Range("I13").Value = 1.036570147524426
$ws.Range("J13").Value = 1.027334242077681
$ws.Range("K13").Value = 1.045987628288445
$ws.Range("L13").Value = 1.024901222023235
$ws.Range("M13").Value = 1.050527817460275
$ws.Range("N13").Value = 1.028793174538316
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.021087261645565
$ws.Range("D14").Value = 1.042567949611854
$ws.Range("E14").Value = 1.021675886391914
$ws.Range("F14").Value = 1.047217768083049
$ws.Range("I14").Value = 1.036615699377189
$ws.Range("J14").Value = 1.027764362496614
$ws.Range("K14").Value = 1.046114568205916
$ws.Range("L14").Value = 1.025300697783504
$ws.Range("M14").Value = 1.050747452853847
$ws.Range("N14").Value = 1.02922390577758
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.021391713026634
$ws.Range("D15").Value = 1.042665767165422
$ws.Range("E15").Value = 1.021942555682682
$ws.Range("F15").Value = 1.047372741433987
$ws.Range("I15").Value = 1.036643692754365
$ws.Range("J15").Value = 1.028029072783126
$ws.Range("K15").Value = 1.04619267152194
$ws.Range("L15").Value = 1.025546539141887
$ws.Range("M15").Value = 1.050882607932811
$ws.Range("N15").Value = 1.02948899198306
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.023159000406712
$ws.Range("D16").Value = 1.043233692118843
$ws.Range("E16").Value = 1.023490387874556
$ws.Range("F16").Value = 1.04827236550514
$ws.Range("I16").Value = 1.036805522764258
$ws.Range("J16").Value = 1.029565261945846
$ws.Range("K16").Value = 1.046645621085476
$ws.Range("L16").Value = 1.026973086242703
$ws.Range("M16").Value = 1.051666703805601
$ws.Range("N16").Value = 1.031027362710595
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.024263460143914
$ws.Range("D17").Value = 1.043588702904453
$ws.Range("E17").Value = 1.024457577890755
$ws.Range("F17").Value = 1.048834599123327
$ws.Range("I17").Value = 1.036906068217264
$ws.Range("J17").Value = 1.030524937590616
$ws.Range("K17").Value = 1.046928308034109
$ws.Range("L17").Value = 1.027864143534453
$ws.Range("M17").Value = 1.052156312454499
$ws.Range("N17").Value = 1.031988401204857
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.024906218617781
$ws.Range("D18").Value = 1.043795335283112
$ws.Range("E18").Value = 1.025020405370191
$ws.Range("F18").Value = 1.049161802306131
$ws.Range("I18").Value = 1.036964368693406
$ws.Range("J18").Value = 1.031083306433438
$ws.Range("K18").Value = 1.047092681961315
$ws.Range("L18").Value = 1.028382543474883
$ws.Range("M18").Value = 1.052441097670393
$ws.Range("N18").Value = 1.032547562995481
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.025125138608876
$ws.Range("D19").Value = 1.043865717533195
$ws.Range("E19").Value = 1.025212093672988
$ws.Range("F19").Value = 1.049273245887941
$ws.Range("I19").Value = 1.036984189172308
$ws.Range("J19").Value = 1.031273461666548
$ws.Range("K19").Value = 1.047148642750454
$ws.Range("L19").Value = 1.028559079524987
$ws.Range("M19").Value = 1.052538068107
$ws.Range("N19").Value = 1.032737988270834
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.024145113107971
$ws.Range("D20").Value = 1.043550659171643
$ws.Range("E20").Value = 1.024353944429347
$ws.Range("F20").Value = 1.048774353359317
$ws.Range("I20").Value = 1.036895316478472
$ws.Range("J20").Value = 1.030422118182203
$ws.Range("K20").Value = 1.046898031531966
$ws.Range("L20").Value = 1.027768680498212
$ws.Range("M20").Value = 1.052103864554372
$ws.Range("N20").Value = 1.031885435781091
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.020941667004356
$ws.Range("D21").Value = 1.042521173434342
$ws.Range("E21").Value = 1.021548357493826
$ws.Range("F21").Value = 1.047143657464637
$ws.Range("I21").Value = 1.036602300633713
$ws.Range("J21").Value = 1.027637765619774
$ws.Range("K21").Value = 1.046077210207956
$ws.Range("L21").Value = 1.025183122431958
$ws.Range("M21").Value = 1.050682811079494
$ws.Range("N21").Value = 1.029097129118653
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.018913826699085
$ws.Range("D22").Value = 1.041869825165271
$ws.Range("E22").Value = 1.019771979213502
$ws.Range("F22").Value = 1.046111502129803
$ws.Range("I22").Value = 1.036414904716822
$ws.Range("J22").Value = 1.02587404793262
$ws.Range("K22").Value = 1.045556400215973
$ws.Range("L22").Value = 1.023544925780072
$ws.Range("M22").Value = 1.049781961010995
$ws.Range("N22").Value = 1.027330906750083
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.019990152253725
$ws.Range("D23").Value = 1.042215508445552
$ws.Range("E23").Value = 1.020714871378953
$ws.Range("F23").Value = 1.046659329523128
$ws.Range("I23").Value = 1.036514549242649
$ws.Range("J23").Value = 1.026810294352715
$ws.Range("K23").Value = 1.045832944894269
$ws.Range("L23").Value = 1.024414579046485
$ws.Range("M23").Value = 1.050260228460999
$ws.Range("N23").Value = 1.028268482747483
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.024198593528094
$ws.Range("D24").Value = 1.043567850854096
$ws.Range("E24").Value = 1.02440077599861
$ws.Range("F24").Value = 1.048801578106341
$ws.Range("I24").Value = 1.036900175795231
$ws.Range("J24").Value = 1.030468582149568
$ws.Range("K24").Value = 1.046911713763283
$ws.Range("L24").Value = 1.027811820267903
$ws.Range("M24").Value = 1.052127565958656
$ws.Range("N24").Value = 1.031931965732618
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.029025502549337
$ws.Range("D25").Value = 1.045119985192356
$ws.Range("E25").Value = 1.028626608219009
$ws.Range("F25").Value = 1.051258688202265
$ws.Range("I25").Value = 1.03733412474149
$ws.Range("J25").Value = 1.034659403968393
$ws.Range("K25").Value = 1.04814350962181
$ws.Range("L25").Value = 1.031701853411071
$ws.Range("M25").Value = 1.054263414790475
$ws.Range("N25").Value = 1.036128738999121
